# Applies the "too many, but highlight new screen size sys" commit:
#  - Admin sheet: remove the "User Table Main Display Edit Capabilities" row
#    (row 5), which shifts every following row up by one.
#  - Scouting Admin sheet: bump the "Manage Users Modal" completion date to
#    the next day, then append the new "Scout Field Schedule" rows (table,
#    modal, Save/Edit/Copy/Delete/Notify).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Admin"
# ---------------------------------------------------------------------
$admin = $wb.Worksheets.Item("Admin")

# Delete the "User Table Main Display Edit Capabilities" row; everything
# below (Manage Users Modal, Error Log, ... Team Contact Form items) shifts
# up by one row automatically, carrying its styles/merges along.
$admin.Rows(5).Delete()

# Match the author's final selection (whole row 5, "Manage Users Modal").
$admin.Activate()
$admin.Rows(5).Select()

# ---------------------------------------------------------------------
# Sheet 2: "Scouting Admin"
# ---------------------------------------------------------------------
$scouting = $wb.Worksheets.Item("Scouting Admin")

# Stamp the date-column formatting (style of B3, "m/d/yyyy") onto the rows
# we're about to fill in (B4 through the new B10) before writing values.
$scouting.Range("B3").Copy()
$scouting.Range("B4:B10").PasteSpecial(-4122)
$scouting.Application.CutCopyMode = $false

# "Manage Users Modal" work date moves from 1/23/2024 to 1/24/2024.
$scouting.Range("B4").Value = 45315

# New "Scout Field Schedule" feature rows, all logged the same day.
$scouting.Range("A5").Value = "Scout Field Schedule Table "
$scouting.Range("B5").Value = 45315

$scouting.Range("A6").Value = "Scout Field Schedule Modal"
$scouting.Range("B6").Value = 45315

$scouting.Range("A7").Value = "     Save"
$scouting.Range("B7").Value = 45315

$scouting.Range("A8").Value = "     Edit"
$scouting.Range("B8").Value = 45315

$scouting.Range("A9").Value = "     Copy"
$scouting.Range("B9").Value = 45315

$scouting.Range("A10").Value = "     Delete"
$scouting.Range("B10").Value = 45315

$scouting.Range("A11").Value = "     Notify"

# Match the author's final selection (single cell A11, "     Notify").
$scouting.Activate()
$scouting.Range("A11").Select()
